$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.386.57"
$ws.Range("E2").Value = "  +2.44%  "

$ws.Range("D3").Value = "2.517.46"
$ws.Range("E3").Value = "  +2.20%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.58%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  +0.32%  "

$ws.Range("D9").Value = "2.521.09"
$ws.Range("E9").Value = "  +2.37%  "

$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("E11").Value = "  -1.45%  "

$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "2.979.27"
$ws.Range("E16").Value = "  +2.42%  "

$ws.Range("D17").Value = "64.158.61"
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("D18").Value = "2.527.45"
$ws.Range("E18").Value = "  +2.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "328.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "641.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.98%  "

$ws.Range("D29").Value = "2.653.25"
$ws.Range("E29").Value = "  +2.38%  "

$ws.Range("E30").Value = "  +4.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.136"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.78%  "

$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "163.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.43%  "

$ws.Range("E46").Value = "  -3.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.73%  "

$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("E49").Value = "  +3.48%  "

$ws.Range("E50").Value = "  +1.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0517"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.24%  "
